$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update adapter-driver label and metrics
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 24
$ws.Range("D5").Value = 97.59999999999999

# Row 6: update adapter-driver label and metrics
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 98.8

# Row 7: update totals (Critical Minutes)
$ws.Range("C7").Value = 73
